$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bank0x3")

# Remove the old "Message Buffer" (row 3) and "Interpolation Buffer" (row 4) rows.
# Excel shifts everything below up by two rows (rows 8/9/10/12 -> 6/7/8/10) and
# the SUM formula on the totals row automatically re-ranges itself.
$ws.Range("A3:I4").EntireRow.Delete()

# The "code" bank (row 2) now covers the whole 0x13 bank (8192 bytes) and is
# relabeled to reflect that it also tracks the text buffers.
$ws.Range("C2").Value = 8192
$ws.Range("A2").Value = "code/buffers"

# Re-purpose the old "Interpolation Buffer" text as a new explanatory note
# row at the bottom of the table, wrapped and taller.
$ws.Range("A13").Value = "This bank is for text processing and includes the text processing buffers. The address of the buffers (textbuffer1,textbuffer2 and interpolation buffer) and tracked by the compiler"
$ws.Range("A13").WrapText = $true
$ws.Rows.Item(13).RowHeight = 145

# Match the author's final selection.
$ws.Activate()
$ws.Range("F4").Select()
